# Updates to Malaysia & Sabah
# Extend the active-cases table (the workbook's single worksheet) from
# row 252 down through row 261 with new daily figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data, rows 253-261:
#   A = Date (serial)
#   B = New Cases
#   C = Imported cases
#   D = Active cases hospitalised
#   F = Total ICU usage
#   G = Ventilator usage
# Columns E (Active cases excl ICU = D-F) and H (% ICU usage = F/D) are
# formulas, same pattern as the rest of the table (E77.. / H24..).
$data = @(
    @(253, 44130, 1240, 2, 9744, 94, 31),
    @(254, 44131,  835, 5, 9903, 89, 32),
    @(255, 44132,  801, 2, 10123, 94, 25),
    @(256, 44133,  649, 7, 10087, 106, 23),
    @(257, 44134,  799, 0, 10392, 90, 20),
    @(258, 44135,  659, 1, 10051, 83, 19),
    @(259, 44136,  957, 0, 10036, 97, 27),
    @(260, 44137,  834, 3, 9968, 91, 32),
    @(261, 44138, 1054, 14, 10135, 94, 32)
)

# Carry the existing row 252 formatting (borders/fonts/number formats)
# down into the new rows before filling in the actual data.
$ws.Range("A252:H252").Copy($ws.Range("A253:H261"))

foreach ($r in $data) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
    $ws.Range("E$row").Formula = "=D$row-F$row"
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Range("H$row").Formula = "=F$row/D$row"
}
